# Preparation for transport:
#  - Independent num/denum conversion
#  - Added some passenger convs
#  - CAP2ACT is now entity dependent
#
# Concretely: insert a new data row (capacity_to_activity) right after the
# existing "buildrate" row (old row 9), before the old "co2_factor" row
# (old row 10). Everything from the old row 10 onward shifts down by one
# row. Update the AutoFilter / _FilterDatabase defined name ranges and the
# selected cell to reflect the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 - this shifts rows 10..429 down to 11..430
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row with the capacity_to_activity parameter
$ws.Range("A10").Value = "CHE"
$ws.Range("B10").Value = "conv_chp_waste"
$ws.Range("C10").Value = "capacity_to_activity"
$ws.Range("D10").Value = "constant"
$ws.Range("G10").Value = 0.001
$ws.Range("H10").Value = "GW/TWh"

# Refresh the AutoFilter so that its range grows from L849 to L850
$ws.AutoFilterMode = $false
[void]$ws.Range("A5:L850").AutoFilter()

# Update the _FilterDatabase defined name to match the new range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$5:`$L`$850"
    }
}

# Update the selected cell shown in the saved view
[void]$ws.Range("C15").Select()
